$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update worker identification fields (row 16)
$ws.Range("C16").Value = "30840363"
$ws.Range("D16").Value = "CARMEN ALICIA MORALES CARRASQUILLA"
$ws.Range("E16").Value = "2508"

# Update monetary values
$ws.Range("E11").Value = 160000
$ws.Range("F16").Value = 160000
$ws.Range("G16").Value = 4000000
